# Onboarding sheet: add "Gross" column before the EPF Employer Contribution
# block, and a "CTC" column before the ESIC Employee column. Both new
# columns inherit the formatting/width of their left neighbour, matching
# what Excel's own "Insert Column" does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "Gross" at column Q (17) ---------------------------------
$leftWidth = $ws.Columns.Item(16).ColumnWidth
$ws.Columns.Item(17).EntireColumn.Insert()
$ws.Range("Q1").Value = "Gross"
$ws.Columns.Item(17).ColumnWidth = $leftWidth

# --- Insert "CTC" at column V (22, after the first insert) -----------
$leftWidth2 = $ws.Columns.Item(21).ColumnWidth
$ws.Columns.Item(22).EntireColumn.Insert()
$ws.Range("V1").Value = "CTC"
$ws.Columns.Item(22).ColumnWidth = $leftWidth2

# --- Restore the view: scrolled right, I17 selected -------------------
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("I17").Select() | Out-Null
